$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- Shapes whose position (off x/y) changes, size unchanged ---
$sh = $s.Shapes.Item("Abgerundetes Rechteck 126")
$sh.Left = 483.87401574803147
$sh.Top = 114.06976777952757

$sh = $s.Shapes.Item("Abgerundetes Rechteck 127")
$sh.Left = 483.87401574803147
$sh.Top = 163.69370278740158

$sh = $s.Shapes.Item("Abgerundetes Rechteck 129")
$sh.Left = 348.66142332283465
$sh.Top = 16.873229346456693

$sh = $s.Shapes.Item("Abgerundetes Rechteck 130")
$sh.Left = 348.66142332283465
$sh.Top = 114.06976777952757

$sh = $s.Shapes.Item("Abgerundetes Rechteck 132")
$sh.Left = 348.66142332283465
$sh.Top = 162.46299212598424

$sh = $s.Shapes.Item("Abgerundetes Rechteck 138")
$sh.Left = 348.66142332283465
$sh.Top = 38.75803249606299

$sh = $s.Shapes.Item("Abgerundetes Rechteck 141")
$sh.Left = 348.66142332283465
$sh.Top = 136.79905511811023

$sh = $s.Shapes.Item("Abgerundetes Rechteck 142")
$sh.Left = 483.87401574803147
$sh.Top = 136.79905511811023

$sh = $s.Shapes.Item("Abgerundetes Rechteck 143")
$sh.Left = 483.87401574803147
$sh.Top = 191.2035453070866

$sh = $s.Shapes.Item("Abgerundetes Rechteck 144")
$sh.Left = 348.66142332283465
$sh.Top = 191.11669291338583

$sh = $s.Shapes.Item("Abgerundetes Rechteck 61")
$sh.Left = 483.87401574803147
$sh.Top = 89.99874015748031

$sh = $s.Shapes.Item("Abgerundetes Rechteck 62")
$sh.Left = 348.66142332283465
$sh.Top = 89.99874015748031

# --- Shapes whose size (ext cx/cy) changes, position unchanged ---
$sh = $s.Shapes.Item("Gerade Verbindung 331")
$sh.Width = 155.5724419448819
$sh.Height = 0.0

$sh = $s.Shapes.Item("Gerade Verbindung 351")
$sh.Width = 56.161890763779525
$sh.Height = 0.0

$sh = $s.Shapes.Item("Gerade Verbindung 355")
$sh.Width = 56.161890763779525
$sh.Height = 22.729291338582676

$sh = $s.Shapes.Item("Gerade Verbindung 374")
$sh.Width = 56.161890763779525
$sh.Height = 0.0

$sh = $s.Shapes.Item("Gerade Verbindung 164")
$sh.Width = 56.24204924409449
$sh.Height = 0.06637895275590551

$sh = $s.Shapes.Item("Gerade Verbindung 183")
$sh.Width = 155.5724419448819
$sh.Height = 21.8848031496063

$sh = $s.Shapes.Item("Gerade Verbindung 65")
$sh.Width = 56.161890763779525
$sh.Height = 24.071023622047242

# --- Connector shapes: full transform change (position, size, rotation/flip) ---
$sh = $s.Shapes.Item("Gerade Verbindung 146")
$sh.Rotation = 0
$sh.HorizontalFlip = -1
$sh.VerticalFlip = 0
$sh.Left = 596.3748031496064
$sh.Top = 99.2603957007874
$sh.Width = 46.10787401574803
$sh.Height = 23.246929133858266

$sh = $s.Shapes.Item("Gerade Verbindung 147")
$sh.Rotation = 0
$sh.HorizontalFlip = -1
$sh.VerticalFlip = 0
$sh.Left = 618.8749696299212
$sh.Top = 114.97779527559055
$sh.Width = 49.85504037007874
$sh.Height = 57.15346456692913

$sh = $s.Shapes.Item("Gerade Verbindung 149")
$sh.Rotation = 0
$sh.HorizontalFlip = -1
$sh.VerticalFlip = -1
$sh.Left = 453.28740157480314
$sh.Top = 25.310788401574804
$sh.Width = 189.19527559055118
$sh.Height = 73.9496062992126

$sh = $s.Shapes.Item("Gerade Verbindung 150")
$sh.Rotation = 0
$sh.HorizontalFlip = -1
$sh.VerticalFlip = -1
$sh.Left = 453.28740157480314
$sh.Top = 47.195590551181105
$sh.Width = 189.19527559055118
$sh.Height = 52.0648031496063

$sh = $s.Shapes.Item("Gerade Verbindung 152")
$sh.Rotation = 0
$sh.HorizontalFlip = -1
$sh.VerticalFlip = 0
$sh.Left = 618.8749696299212
$sh.Top = 114.97779527559055
$sh.Width = 49.85504037007874
$sh.Height = 84.66330708661417

$sh = $s.Shapes.Item("Gerade Verbindung 153")
$sh.Rotation = 0
$sh.HorizontalFlip = -1
$sh.VerticalFlip = 0
$sh.Left = 618.8749696299212
$sh.Top = 99.2603957007874
$sh.Width = 23.60771753543307
$sh.Height = 45.97622147244094

$sh = $s.Shapes.Item("Gerade Verbindung 68")
$sh.Rotation = 0
$sh.HorizontalFlip = -1
$sh.VerticalFlip = -1
$sh.Left = 596.3748031496064
$sh.Top = 98.43629921259843
$sh.Width = 46.10787401574803
$sh.Height = 0.8240954881889764
